# Adds the 2023-24 UEFA Champions League final (Real Madrid beat Borussia
# Dortmund 2-0 at Wembley Stadium, London) as a new row, and freezes the
# header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count   # 69 -> header + 68 finals
$newRow = $lastRow + 1

$ws.Range("A$newRow").Value = "2023-24"
$ws.Range("B$newRow").Value = "Spain"
$ws.Range("C$newRow").Value = "Real Madrid"
$ws.Range("D$newRow").Value = 2
$ws.Range("E$newRow").Value = 0
$ws.Range("F$newRow").Value = "Borussia Dortmund"
$ws.Range("G$newRow").Value = "Germany"
$ws.Range("H$newRow").Value = "Wembley Stadium"
$ws.Range("I$newRow").Value = "London"
$ws.Range("J$newRow").Value = "England"
$ws.Range("K$newRow").Value = 86212

# Columns L:N hold the text "True"/"False" (not real booleans) throughout
# this sheet, and this final (like the one before it) was decided in
# normal time: normal-time=True, extra-time=False, penalty=False -- the
# same text already sitting in row $lastRow. Assigning the bare word via
# .Value/.Formula gets auto-coerced into a real Boolean, so copy the text
# straight out of the previous row instead, which preserves its text type.
$srcRange = "L" + $lastRow + ":N" + $lastRow
$ws.Range($srcRange).Copy() | Out-Null
$ws.Range("L$newRow").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

# Freeze the header row, same as the saved view in the target workbook.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
